$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '89.552.06'
$ws.Range('E2').Value = '  -1.53%  '
$ws.Range('D3').Value = '3.103.66'
$ws.Range('E3').Value = '  -2.03%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.997'
$ws.Range('D4').ClearFormats()
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.64'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.27%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '623.03'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.70%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.373'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -4.62%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.824'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +16.01%  '
$ws.Range('E9').Value = '  +0.00%  '
$ws.Range('D10').Value = '3.101.85'
$ws.Range('E10').Value = '  -1.91%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.611'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +8.99%  '
$ws.Range('E12').Value = '  +1.19%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000243'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -3.08%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.32'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.74%  '
$ws.Range('D15').Value = '89.197.64'
$ws.Range('E15').Value = '  -1.38%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '32.35'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.23%  '
$ws.Range('D17').Value = '3.665.29'
$ws.Range('E17').Value = '  -2.34%  '
$ws.Range('D18').Value = '3.117.16'
$ws.Range('E18').Value = '  -1.78%  '
$ws.Range('E19').Value = '  +3.80%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0000216'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +4.08%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.51'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +2.47%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '426.44'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.93%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.31'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -1.10%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.95'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.14%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.55'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +7.68%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.07'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +4.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '83.94'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +4.88%  '
$ws.Range('D28').Value = '3.252.80'
$ws.Range('E28').Value = '  -2.63%  '
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.164'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +5.52%  '
$ws.Range('E31').Value = '  +1.38%  '
$ws.Range('E32').Value = '  -1.09%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '512.09'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.70%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.71'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -6.05%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.73'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -2.04%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.27'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.65%  '
$ws.Range('E37').Value = '  -3.44%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '22.46'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +1.25%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '22.29'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.38%  '
$ws.Range('E40').Value = '  +2.37%  '
$ws.Range('E41').Value = '  +0.15%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.365'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.41%  '
$ws.Range('B44').Value = 'Stacks'
$ws.Range('C44').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.84'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -3.04%  '
$ws.Range('B45').Value = 'Stellar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.135'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +7.05%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '145.31'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.91%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0701'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +12.97%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '43.48'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.85%  '
$ws.Range('E49').Value = '  +2.44%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '160.48'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -5.55%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.708'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -3.81%  '
